$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.108.70'
$ws.Range('E2').Value = '  +0.24%  '

$ws.Range('D3').Value = '1.878.62'
$ws.Range('E3').Value = '  -1.04%  '

$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = '313.41'
$ws.Range('E5').Value = '  +0.26%  '

$ws.Range('E6').Value = '  +0.01%  '

$ws.Range('D7').Value = '0.5088'
$ws.Range('E7').Value = '  +0.56%  '

$ws.Range('D8').Value = '0.3848'

$ws.Range('D9').Value = '0.08992'
$ws.Range('E9').Value = '  -2.06%  '

$ws.Range('D10').Value = '1.124'
$ws.Range('E10').Value = '  -0.82%  '

$ws.Range('D11').Value = '41.57'
$ws.Range('E11').Value = '  -0.56%  '

$ws.Range('D12').Value = '6.343'
$ws.Range('E12').Value = '  -0.18%  '

$ws.Range('D13').Value = '20.77'
$ws.Range('E13').Value = '  +0.14%  '

$ws.Range('D14').Value = '1.877.76'
$ws.Range('E14').Value = '  -1.11%  '

$ws.Range('D15').Value = '7.219'
$ws.Range('E15').Value = '  -0.83%  '

$ws.Range('D16').Value = '1.002'
$ws.Range('E16').Value = '  -0.02%  '

$ws.Range('D17').Value = '0.00001108'
$ws.Range('E17').Value = '  -0.55%  '

$ws.Range('D18').Value = '91.25'
$ws.Range('E18').Value = '  -1.05%  '

$ws.Range('D19').Value = '0.06597'
$ws.Range('E19').Value = '  +0.28%  '

$ws.Range('D20').Value = '18.18'
$ws.Range('E20').Value = '  +2.47%  '

$ws.Range('E21').Value = '  +0.01%  '

$ws.Range('D22').Value = '6.120'
$ws.Range('E22').Value = '  -1.47%  '

$ws.Range('D23').Value = '28.131.56'
$ws.Range('E23').Value = '  +0.14%  '

$ws.Range('D24').Value = '11.41'
$ws.Range('E24').Value = '  +0.87%  '

$ws.Range('D25').Value = '2.275'
$ws.Range('E25').Value = '  -1.95%  '

$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').Value = '2.094.05'
$ws.Range('E26').Value = '  -1.05%  '

$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '2.548'
$ws.Range('E27').Value = '  -2.19%  '

$ws.Range('D28').Value = '20.79'

$ws.Range('D29').Value = '156.92'
$ws.Range('E29').Value = '  -0.16%  '

$ws.Range('D30').Value = '126.91'
$ws.Range('E30').Value = '  -0.12%  '

$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '1.063'
$ws.Range('E31').Value = '  -1.71%  '

$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '0.1052'
$ws.Range('E32').Value = '  -1.31%  '

$ws.Range('D33').Value = '5.612'
$ws.Range('E33').Value = '  +0.18%  '

$ws.Range('D34').Value = '3.600'
$ws.Range('E34').Value = '  -0.42%  '

$ws.Range('D35').Value = '9.643'
$ws.Range('E35').Value = '  +0.77%  '

$ws.Range('D36').Value = '0.06591'
$ws.Range('E36').Value = '  -0.73%  '

$ws.Range('D37').Value = '0.02418'
$ws.Range('E37').Value = '  +0.66%  '

$ws.Range('D38').Value = '0.2182'
$ws.Range('E38').Value = '  +0.82%  '

$ws.Range('D39').Value = '1.277'
$ws.Range('E39').Value = '  +1.69%  '

$ws.Range('D40').Value = '1.211'
$ws.Range('E40').Value = '  -0.76%  '

$ws.Range('D41').Value = '0.6403'
$ws.Range('E41').Value = '  +0.65%  '

$ws.Range('D42').Value = '11.50'
$ws.Range('E42').Value = '  +1.02%  '

$ws.Range('D43').Value = '4.924'
$ws.Range('E43').Value = '  -1.25%  '

$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  -0.01%  '

$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.6039'
$ws.Range('E45').Value = '  +1.08%  '

$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '13.19'
$ws.Range('E46').Value = '  -0.71%  '

$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').Value = '3.675'
$ws.Range('E47').Value = '  -0.81%  '

$ws.Range('B48').Value = 'WEMIXTOKEN'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '1.274'
$ws.Range('E48').Value = '  +0.01%  '

$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').Value = '1.243'
$ws.Range('E49').Value = '  +5.63%  '

$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = '2.001'
$ws.Range('E50').Value = '  -0.43%  '

$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').Value = '121.42'
$ws.Range('E51').Value = '  -0.76%  '
